$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.987.16"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "2.295.79"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "310.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.14%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("E14").Value = "  +21.61%  "
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "2.634.72"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "2.288.63"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "43.023.42"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "257.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0903"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.52%  "
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.46%  "
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.07%  "
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.82%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0985"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.53%  "
